$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert three new columns (J, K, L) before the old "comment" column, which
# shifts the old J column (comments) to M, carrying its data/width along.
$ws.Columns("J:L").Insert()

# New header cells for the inserted columns.
$ws.Range("J1").Value = "Mutationstrength"
$ws.Range("K1").Value = "Weight decay fac"
$ws.Range("L1").Value = "Create up to"

# Fill the new J/K/L values for the existing data rows (2-28).
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 10).Value = 1.6
    $ws.Cells.Item($r, 11).Value = 0
    $ws.Cells.Item($r, 12).Value = 250
}

# Append three new rows of data (29-31). Copy an existing similarly
# formatted row (A:I) down as a template, so number formats / text vs.
# boolean typing (e.g. the "true" flag) are preserved exactly, then tweak
# the few cells that differ.
$ws.Range("A24:I24").Copy($ws.Range("A29:I29"))
$ws.Range("A24:I24").Copy($ws.Range("A30:I30"))
$ws.Range("A24:I24").Copy($ws.Range("A31:I31"))

for ($r = 29; $r -le 31; $r++) {
    $ws.Cells.Item($r, 6).ClearContents() | Out-Null
    $ws.Cells.Item($r, 7).Value = 10
    $ws.Cells.Item($r, 8).Value = 10
}

$ws.Cells.Item(29, 10).Value = 0.001
$ws.Cells.Item(29, 11).Value = 0.00005
$ws.Cells.Item(29, 12).Value = 580
$ws.Cells.Item(29, 13).Value = "nach ca 52 M Comp. Mit minimaler neuron count => 1.0 rating"

$ws.Cells.Item(30, 10).Value = 0.001
$ws.Cells.Item(30, 11).Value = 0.000005
$ws.Cells.Item(30, 12).Value = 580
$ws.Cells.Item(30, 13).Value = "nach ca 70 M Comp. Mit minimaler neuron count => 1.0 rating"

$ws.Cells.Item(31, 10).Value = 0.001
$ws.Cells.Item(31, 11).Value = 0.0000005
$ws.Cells.Item(31, 12).Value = 580
$ws.Cells.Item(31, 13).Value = "nach ca 47 M Comp. Mit minimaler neuron count => 1.0 rating"

# Column widths: bestFit/autofit per the target layout (values chosen so the
# resulting stored width is as close as possible to the recorded bestFit
# widths of 15.21875 / 15 / 47.33203125 given the engine's internal
# character-width quantization).
$ws.Columns("J").ColumnWidth = 14.25
$ws.Columns("K").ColumnWidth = 14.1
$ws.Columns("M").ColumnWidth = 46.42

# View tweaks recorded in the saved file.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("M31").Select() | Out-Null
